$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - 2021-10-28, Code, new note about planning algorithm
$ws.Range("A8").Value = 44497
$ws.Range("B8").Value = 0.375
$ws.Range("C8").Value = 0.46875
$ws.Range("E8").Value = "Code"
$ws.Range("F8").Value = "Writing code from the UML diagrams. Specifically working on planning algorithm."

# Row 9 - 2021-10-29, Code, new note about Vision doc / team velocity
$ws.Range("A9").Value = 44498
$ws.Range("B9").Value = 0.46875
$ws.Range("C9").Value = 0.5
$ws.Range("E9").Value = "Code"
$ws.Range("F9").Value = "Writing code from the UML diagrams. Specifically working on planning algorithm. Also updated Vision doc to add team velocity input."

# Row 10 - 2021-10-30, Code, reuse of existing "Writing code..." note
$ws.Range("A10").Value = 44499
$ws.Range("B10").Value = 0.625
$ws.Range("C10").Value = 0.77083333333333337
$ws.Range("E10").Value = "Code"
$ws.Range("F10").Value = "Writing code from the UML diagrams."

# Row 11 - 2021-11-01, Code, reuse of existing "Writing code..." note
$ws.Range("A11").Value = 44501
$ws.Range("B11").Value = 0.33333333333333331
$ws.Range("C11").Value = 0.70833333333333337
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("D11").NumberFormat = "h:mm:ss;@"
$ws.Range("E11").Value = "Code"
$ws.Range("F11").Value = "Writing code from the UML diagrams."

# Selection moves to B9, matching the saved cursor position in the workbook
$ws.Range("B9").Select() | Out-Null
